$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update currency quote values (column B)
$ws.Range("B2").Value = 5.23
$ws.Range("B3").Value = 3.39
$ws.Range("B4").Value = 3.83
$ws.Range("B5").Value = 5.62
$ws.Range("B6").Value = 5.79
$ws.Range("B7").Value = 0.0339
$ws.Range("B8").Value = 6.52
$ws.Range("B9").Value = 1.39
$ws.Range("B15").Value = 0.72

# Update report date/time
$ws.Range("C18").Value = "18/04/2024"
$ws.Range("D18").Value = "03:18"
